$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data-row formatting: the placeholder row (A2:F2) was bold+centered by
# mistake; normalize it to a plain, left-aligned data row. ---
$ws.Range("A2:F2").Font.Bold = $false
$ws.Range("A2:F2").HorizontalAlignment = -4131   # xlLeft

# --- Header row: keep bold, add a thin box border around each header cell. ---
$ws.Range("A1:F1").Borders.LineStyle = 1         # xlContinuous (defaults to thin)

# --- New "Total" label + SUM formula, integrated into the report. ---
$ws.Range("H3").Value = "Total"
$ws.Range("H3").Font.Bold = $true
$ws.Range("H3").HorizontalAlignment = -4108      # xlCenter
$ws.Range("H3").Borders.LineStyle = 1

$ws.Range("H4").Formula = "=SUM(F2:F1000)"
$ws.Range("H4").Borders.LineStyle = 1

# --- Keep the active selection parked below the new total, like the source file. ---
$ws.Range("H5").Select()
